$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price"): cells whose new value is a plain decimal number must be
# pre-formatted as Text, otherwise Excel silently parses the literal string into a
# float (dropping trailing zeros / switching to scientific notation for tiny values).
# NumberFormat is reset back to General (via Style) once the literal text is stored,
# so the cell ends up with no explicit style - matching the source file.

$ws.Range("D2").Value = '26.793.92'
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("D3").Value = '1.868.35'
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("E4").Value = '  -0.11%  '
$d = $ws.Range("D5")
$d.NumberFormat = '@'
$d.Value = '300.10'
$d.Style = 'Normal'
$ws.Range("E5").Value = '  -2.44%  '
$ws.Range("E6").Value = '  -0.08%  '
$d = $ws.Range("D7")
$d.NumberFormat = '@'
$d.Value = '0.5370'
$d.Style = 'Normal'
$ws.Range("E7").Value = '  +1.94%  '
$d = $ws.Range("D8")
$d.NumberFormat = '@'
$d.Value = '0.3737'
$d.Style = 'Normal'
$ws.Range("E9").Value = '  -2.20%  '
$d = $ws.Range("D10")
$d.NumberFormat = '@'
$d.Value = '21.54'
$d.Style = 'Normal'
$ws.Range("E10").Value = '  -2.56%  '
$d = $ws.Range("D11")
$d.NumberFormat = '@'
$d.Value = '0.8866'
$d.Style = 'Normal'
$ws.Range("E11").Value = '  -1.69%  '
$d = $ws.Range("D12")
$d.NumberFormat = '@'
$d.Value = '0.08151'
$d.Style = 'Normal'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").Value = '1.894.82'
$ws.Range("E13").Value = '  +42.33%  '
$d = $ws.Range("D14")
$d.NumberFormat = '@'
$d.Value = '92.34'
$d.Style = 'Normal'
$ws.Range("E14").Value = '  -3.61%  '
$d = $ws.Range("D15")
$d.NumberFormat = '@'
$d.Value = '5.288'
$d.Style = 'Normal'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("E16").Value = '  -0.13%  '
$d = $ws.Range("D17")
$d.NumberFormat = '@'
$d.Value = '14.81'
$d.Style = 'Normal'
$ws.Range("E17").Value = '  +0.28%  '
$d = $ws.Range("D18")
$d.NumberFormat = '@'
$d.Value = '0.000008485'
$d.Style = 'Normal'
$ws.Range("E18").Value = '  -1.80%  '
$d = $ws.Range("D19")
$d.NumberFormat = '@'
$d.Value = '1.000'
$d.Style = 'Normal'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '26.823.38'
$ws.Range("E20").Value = '  -1.79%  '
$d = $ws.Range("D21")
$d.NumberFormat = '@'
$d.Value = '4.965'
$d.Style = 'Normal'
$ws.Range("E22").Value = '  -1.74%  '
$d = $ws.Range("D23")
$d.NumberFormat = '@'
$d.Value = '6.370'
$d.Style = 'Normal'
$ws.Range("E23").Value = '  -2.26%  '
$d = $ws.Range("D24")
$d.NumberFormat = '@'
$d.Value = '2.281'
$d.Style = 'Normal'
$ws.Range("E24").Value = '  -0.66%  '
$d = $ws.Range("D25")
$d.NumberFormat = '@'
$d.Value = '146.15'
$d.Style = 'Normal'
$d = $ws.Range("D26")
$d.NumberFormat = '@'
$d.Value = '1.737'
$d.Style = 'Normal'
$ws.Range("E26").Value = '  +0.01%  '
$d = $ws.Range("D27")
$d.NumberFormat = '@'
$d.Value = '17.96'
$d.Style = 'Normal'
$ws.Range("E27").Value = '  -1.53%  '
$d = $ws.Range("D28")
$d.NumberFormat = '@'
$d.Value = '113.67'
$d.Style = 'Normal'
$ws.Range("E28").Value = '  -2.28%  '
$d = $ws.Range("D29")
$d.NumberFormat = '@'
$d.Value = '4.696'
$d.Style = 'Normal'
$ws.Range("E29").Value = '  -2.88%  '
$d = $ws.Range("D30")
$d.NumberFormat = '@'
$d.Value = '4.624'
$d.Style = 'Normal'
$ws.Range("E30").Value = '  -4.09%  '
$d = $ws.Range("D31")
$d.NumberFormat = '@'
$d.Value = '0.09102'
$d.Style = 'Normal'
$ws.Range("E31").Value = '  -1.70%  '
$d = $ws.Range("D32")
$d.NumberFormat = '@'
$d.Value = '0.8125'
$d.Style = 'Normal'
$ws.Range("E32").Value = '  -3.07%  '
$d = $ws.Range("D33")
$d.NumberFormat = '@'
$d.Value = '0.05018'
$d.Style = 'Normal'
$d = $ws.Range("D34")
$d.NumberFormat = '@'
$d.Value = '1.169'
$d.Style = 'Normal'
$ws.Range("E34").Value = '  -4.66%  '
$d = $ws.Range("D35")
$d.NumberFormat = '@'
$d.Value = '2.943'
$d.Style = 'Normal'
$ws.Range("E35").Value = '  -2.23%  '
$d = $ws.Range("D36")
$d.NumberFormat = '@'
$d.Value = '0.6072'
$d.Style = 'Normal'
$ws.Range("E36").Value = '  +5.54%  '
$d = $ws.Range("D37")
$d.NumberFormat = '@'
$d.Value = '2.667'
$d.Style = 'Normal'
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  -4.96%  '
$d = $ws.Range("D39")
$d.NumberFormat = '@'
$d.Value = '0.01943'
$d.Style = 'Normal'
$ws.Range("E39").Value = '  -3.07%  '
$ws.Range("E40").Value = '  -1.06%  '
$d = $ws.Range("D41")
$d.NumberFormat = '@'
$d.Value = '0.5282'
$d.Style = 'Normal'
$ws.Range("E41").Value = '  +7.42%  '
$d = $ws.Range("D42")
$d.NumberFormat = '@'
$d.Value = '8.754'
$d.Style = 'Normal'
$ws.Range("E42").Value = '  -6.04%  '
$d = $ws.Range("D43")
$d.NumberFormat = '@'
$d.Value = '6.467'
$d.Style = 'Normal'
$ws.Range("E43").Value = '  -0.95%  '
$d = $ws.Range("D44")
$d.NumberFormat = '@'
$d.Value = '116.26'
$d.Style = 'Normal'
$ws.Range("E44").Value = '  -0.49%  '
$d = $ws.Range("D45")
$d.NumberFormat = '@'
$d.Value = '0.1484'
$d.Style = 'Normal'
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("E46").Value = '  -0.09%  '
$d = $ws.Range("D47")
$d.NumberFormat = '@'
$d.Value = '1.644'
$d.Style = 'Normal'
$ws.Range("E47").Value = '  +0.48%  '
$d = $ws.Range("D48")
$d.NumberFormat = '@'
$d.Value = '9.971'
$d.Style = 'Normal'
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("E49").Value = '  -4.03%  '
$d = $ws.Range("D50")
$d.NumberFormat = '@'
$d.Value = '0.06064'
$d.Style = 'Normal'
$ws.Range("E50").Value = '  -1.87%  '
$d = $ws.Range("D51")
$d.NumberFormat = '@'
$d.Value = '62.03'
$d.Style = 'Normal'
$ws.Range("E51").Value = '  -2.87%  '
